$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new "Red Desert" rows (mirrors the existing Savanna section
# pattern: a "    " indented category row followed by a "        " indented
# detail row). Written in this order so the shared-string table assigns
# index 63 to the 8-space variant and 64 to the 4-space variant.
$ws.Range("A53").Value = "        Red Desert"
$ws.Range("A52").Value = "    Red Desert"

# Row 52 previously had an explicit custom row height (15.75); it should
# revert back to the sheet's default row height with no explicit <row ht=.../>.
$ws.Rows.Item(52).AutoFit()

# Update the active selection/cell from B52 to A55.
$ws.Range("A55").Select() | Out-Null
